$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Nomor Telepon" column (E) for the
# new "Status" column, shifting "Nomor Telepon" and everything after it one
# column to the right.
$ws.Columns.Item(5).Insert() | Out-Null

# New "Status" column (now E)
$ws.Range("E1").Value = "Status (Aktif/Cuti/Resign/Pensiun)"
$ws.Range("E2").Value = "Aktif"

# Newly appended columns G:L with headers + sample row values
$ws.Range("G1").Value = "Agama"
$ws.Range("H1").Value = "Provinsi"
$ws.Range("I1").Value = "Kota/Kabupaten"
$ws.Range("J1").Value = "Kecamatan"
$ws.Range("K1").Value = "Desa/Kelurahan"
$ws.Range("L1").Value = "ID Program Studi"

$ws.Range("G2").Value = "Islam"
$ws.Range("H2").Value = "Jawa Barat"
$ws.Range("I2").Value = "Bandung"
$ws.Range("J2").Value = "Coblong"
$ws.Range("K2").Value = "Dago"
$ws.Range("L2").Value = 1

# Header cells should share the bold/filled style used by the rest of row 1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-fit the (now much longer) "Status" header column and size the new
# columns to their content, mirroring the sheet's bestFit-style columns.
$ws.Range("E1").EntireColumn.ColumnWidth = 40.25
$ws.Range("G1").EntireColumn.ColumnWidth = 6.1
$ws.Range("H1").EntireColumn.ColumnWidth = 11.95
$ws.Range("I1").EntireColumn.ColumnWidth = 16.6
$ws.Range("J1").EntireColumn.ColumnWidth = 10.75
$ws.Range("K1").EntireColumn.ColumnWidth = 16.6
$ws.Range("L1").EntireColumn.ColumnWidth = 19.1

# Match the original workbook's selection/active-cell convention for the
# (now wider) header row.
$ws.Range("A1:L1").Select() | Out-Null
